$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

$updates = @(
    @{Cell='D2'; Value='260.93'}
    @{Cell='G2'; Value='21'}
    @{Cell='G3'; Value='21'}
    @{Cell='D4'; Value='6.188'}
    @{Cell='G4'; Value='21'}
    @{Cell='D5'; Value='0.06110'}
    @{Cell='G5'; Value='21'}
    @{Cell='D6'; Value='6.742'}
    @{Cell='G6'; Value='21'}
    @{Cell='D7'; Value='3.490'}
    @{Cell='G7'; Value='21'}
    @{Cell='G8'; Value='21'}
    @{Cell='D9'; Value='0.7992'}
    @{Cell='G9'; Value='21'}
    @{Cell='G10'; Value='21'}
    @{Cell='D11'; Value='0.08084'}
    @{Cell='G11'; Value='21'}
    @{Cell='D12'; Value='0.03322'}
    @{Cell='G12'; Value='21'}
    @{Cell='D13'; Value='0.03082'}
    @{Cell='G13'; Value='21'}
    @{Cell='D14'; Value='0.09302'}
    @{Cell='G14'; Value='21'}
    @{Cell='D15'; Value='3.923'}
    @{Cell='G15'; Value='21'}
    @{Cell='D16'; Value='0.001708'}
    @{Cell='G16'; Value='21'}
    @{Cell='D17'; Value='0.04828'}
    @{Cell='G17'; Value='21'}
    @{Cell='D18'; Value='0.0006145'}
    @{Cell='G18'; Value='21'}
    @{Cell='D19'; Value='0.006195'}
    @{Cell='G19'; Value='21'}
    @{Cell='D20'; Value='0.001102'}
    @{Cell='G20'; Value='21'}
    @{Cell='D21'; Value='0.003397'}
    @{Cell='G21'; Value='21'}
    @{Cell='G22'; Value='21'}
    @{Cell='D23'; Value='3.693'}
    @{Cell='G23'; Value='21'}
    @{Cell='D24'; Value='2.259'}
    @{Cell='G24'; Value='21'}
    @{Cell='D25'; Value='0.3357'}
    @{Cell='G25'; Value='21'}
    @{Cell='G26'; Value='21'}
    @{Cell='D27'; Value='0.0006171'}
    @{Cell='G27'; Value='21'}
    @{Cell='G28'; Value='21'}
    @{Cell='G29'; Value='21'}
    @{Cell='G30'; Value='21'}
    @{Cell='G31'; Value='21'}
    @{Cell='G32'; Value='21'}
    @{Cell='G33'; Value='21'}
    @{Cell='G34'; Value='21'}
    @{Cell='G35'; Value='21'}
    @{Cell='G36'; Value='21'}
    @{Cell='G37'; Value='21'}
    @{Cell='G38'; Value='21'}
    @{Cell='G39'; Value='21'}
    @{Cell='D40'; Value='0.04591'}
    @{Cell='G40'; Value='21'}
    @{Cell='B41'; Value='KickToken'}
    @{Cell='C41'; Value='https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'}
    @{Cell='D41'; Value='0.007122'}
    @{Cell='E41'; Value='40KickTokenKICK'}
    @{Cell='G41'; Value='21'}
    @{Cell='B42'; Value='BKEXToken'}
    @{Cell='C42'; Value='https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'}
    @{Cell='D42'; Value='0.1119'}
    @{Cell='E42'; Value='41BKEXTokenBKK'}
    @{Cell='G42'; Value='21'}
    @{Cell='B43'; Value='CEJI'}
    @{Cell='C43'; Value='https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'}
    @{Cell='D43'; Value='0.003901'}
    @{Cell='E43'; Value='42CEJICEJI'}
    @{Cell='G43'; Value='21'}
    @{Cell='G44'; Value='21'}
    @{Cell='D45'; Value='0.002973'}
    @{Cell='G45'; Value='21'}
    @{Cell='D46'; Value='0.00006017'}
    @{Cell='G46'; Value='21'}
    @{Cell='D47'; Value='0.00000000750'}
    @{Cell='G47'; Value='21'}
    @{Cell='D48'; Value='0.7507'}
    @{Cell='G48'; Value='21'}
    @{Cell='D49'; Value='0.1154'}
    @{Cell='G49'; Value='21'}
    @{Cell='G50'; Value='21'}
    @{Cell='D51'; Value='0.01010'}
    @{Cell='G51'; Value='21'}
)

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}